$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 3 and Row 4 effectively swap their species-specific data (the Lichenicolous
# "Garnlav" record and the "Granticka" record exchange rows), while the
# "Taxonsorteringsordning" (column B) values are bumped by one for both taxa.
# ---------------------------------------------------------------------------

# Row 3 -> becomes the "Granticka" record (previously in row 4), with B+1
$ws.Range("A3").Value = 130938742
$ws.Range("B3").Value = 91829
$ws.Range("E3").Value = 5432
$ws.Range("F3").Value = "Granticka"
$ws.Range("G3").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H3").Value = ""
$ws.Range("K3").Value = "teleomorf"
$ws.Range("Q3").Value = 476481
$ws.Range("R3").Value = 7033943

# Row 4 -> becomes the "Garnlav" record (previously in row 3), with B+1
$ws.Range("A4").Value = 130938748
$ws.Range("B4").Value = 79244
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = "Garnlav"
$ws.Range("G4").Value = "Alectoria sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("K4").Value = ""
$ws.Range("Q4").Value = 476501
$ws.Range("R4").Value = 7033912

# ---------------------------------------------------------------------------
# Stand-alone "Taxonsorteringsordning" (column B) bumps on otherwise
# untouched rows (same taxon, re-sequenced sort order).
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = 91809
$ws.Range("B7").Value = 79244
$ws.Range("B8").Value = 79244
$ws.Range("B9").Value = 79244
$ws.Range("B11").Value = 79244
$ws.Range("B12").Value = 91829
$ws.Range("B15").Value = 83224
$ws.Range("B17").Value = 91829
$ws.Range("B18").Value = 78256
$ws.Range("B21").Value = 83224
$ws.Range("B24").Value = 79244
$ws.Range("B25").Value = 79244

# ---------------------------------------------------------------------------
# Row 19 and Row 20 swap their species-specific data (the "Garnlav" record
# and the "Tretåig hackspett" record exchange rows), including the extra
# observation detail columns (I-M, AC, AF, AM, AO).
# ---------------------------------------------------------------------------

# Row 19 -> becomes the "Tretåig hackspett" record (previously in row 20)
$ws.Range("A19").Value = 130938734
$ws.Range("B19").Value = 57884
$ws.Range("E19").Value = 100109
$ws.Range("F19").Value = "Tretåig hackspett"
$ws.Range("G19").Value = "Picoides tridactylus"
$ws.Range("H19").Value = "(Linnaeus, 1758)"
$ws.Range("J19").ClearContents()
$ws.Range("L19").Value = ""
$ws.Range("M19").Value = "äldre spår"
$ws.Range("Q19").Value = 476457
$ws.Range("R19").Value = 7033634
$ws.Range("AC19").Value = "Ringhack, äldre, på gran."
$ws.Range("AF19").ClearContents()
$ws.Range("AM19").Value = "Trädstam på levande träd"
$ws.Range("AO19").Value = "Stem on living tree # Picea abies"

# Row 20 -> becomes the "Garnlav" record (previously in row 19)
$ws.Range("A20").Value = 130938752
$ws.Range("B20").Value = 79244
$ws.Range("E20").Value = 6425
$ws.Range("F20").Value = "Garnlav"
$ws.Range("G20").Value = "Alectoria sarmentosa"
$ws.Range("H20").Value = "(Ach.) Ach."
$ws.Range("J20").Value = ""
$ws.Range("L20").ClearContents()
$ws.Range("M20").ClearContents()
$ws.Range("Q20").Value = 476286
$ws.Range("R20").Value = 7033527
$ws.Range("AC20").ClearContents()
$ws.Range("AF20").Value = ""
$ws.Range("AM20").ClearContents()
$ws.Range("AO20").Value = "Picea abies"
